# Update the Views worksheet to reflect the new/edited view rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (id=2): the "columns" cell picks up a value it previously did not have.
$ws.Range("D4").Value = "admited_no, Available, course_id, no_of_seat"

# Row 5 (id=3): view renamed to DEMO1 and its columns list updated.
$ws.Range("B5").Value = "DEMO1"
$ws.Range("D5").Value = "admited_no, alloted_no, no_of_seat"

# Row 6 (id=4): view renamed to FFFF (columns remains blank).
$ws.Range("B6").Value = "FFFF"

# Match the saved selection state from the workbook.
$ws.Range("K11").Select()
